$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5363
$ws.Range("J137").Value = 5074.636
$ws.Range("L137").Value = 15223.908
$ws.Range("N137").Value = -20323.908
$ws.Range("H138").Value = 27785328
$ws.Range("I138").Value = 76926700
$ws.Range("K138").Value = 230780100
$ws.Range("M138").Value = -230774960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 22021.5
$ws.Range("I28").Value = 8724
$ws.Range("K28").Value = 8724
$ws.Range("M28").Value = -8532
$ws.Range("H32").Value = 3614.1746
$ws.Range("I32").Value = 2823.3147
$ws.Range("J32").Value = 8359.333000000001
$ws.Range("K32").Value = 2823.3147
$ws.Range("L32").Value = 8359.333000000001
$ws.Range("M32").Value = -2536.3147
$ws.Range("N32").Value = -8933.333000000001
$ws.Range("H61").Value = 4950.4326
$ws.Range("I61").Value = 4860.853
$ws.Range("K61").Value = 4860.853
$ws.Range("M61").Value = -4648.853
$ws.Range("H74").Value = 2812.0908
$ws.Range("I74").Value = 2122.3333
$ws.Range("J74").Value = 3639.8
$ws.Range("K74").Value = 2122.3333
$ws.Range("L74").Value = 3639.8
$ws.Range("M74").Value = -1248.3333
$ws.Range("N74").Value = -5387.8
$ws.Range("H77").Value = 2812.0908
$ws.Range("I77").Value = 2122.3333
$ws.Range("J77").Value = 3639.8
$ws.Range("K77").Value = 10611.6665
$ws.Range("L77").Value = 18199
$ws.Range("M77").Value = -6243.666499999999
$ws.Range("N77").Value = -26935
$ws.Range("H99").Value = 22021.5
$ws.Range("I99").Value = 8724
$ws.Range("K99").Value = 8724
$ws.Range("M99").Value = -5729
$ws.Range("H132").Value = 10620.786
$ws.Range("I132").Value = 10516
$ws.Range("K132").Value = 31548
$ws.Range("M132").Value = -29018
$ws.Range("H136").Value = 4950.4326
$ws.Range("I136").Value = 4860.853
$ws.Range("K136").Value = 14582.559
$ws.Range("M136").Value = -12032.559

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1227
$ws.Range("I36").Value = 908.75
$ws.Range("K36").Value = 908.75
$ws.Range("M36").Value = -374.75
$ws.Range("H99").Value = 12598.75
$ws.Range("I99").Value = 12598.75
$ws.Range("K99").Value = 12598.75
$ws.Range("M99").Value = -11100.75
$ws.Range("H107").Value = 4599.413
$ws.Range("I107").Value = 3977.1724
$ws.Range("K107").Value = 3977.1724
$ws.Range("M107").Value = -2057.1724
$ws.Range("H133").Value = 179999
$ws.Range("J133").Value = 179999
$ws.Range("L133").Value = 179999
$ws.Range("N133").Value = -190119
$ws.Range("H134").Value = 21331.666
$ws.Range("I134").Value = 21331.666
$ws.Range("K134").Value = 63994.99800000001
$ws.Range("M134").Value = -61459.99800000001
$ws.Range("H135").Value = 144996
$ws.Range("J135").Value = 144996
$ws.Range("L135").Value = 144996
$ws.Range("N135").Value = -155136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1377.5358
$ws.Range("I107").Value = 638.45
$ws.Range("J107").Value = 3225.25
$ws.Range("K107").Value = 638.45
$ws.Range("L107").Value = 3225.25
$ws.Range("M107").Value = 1281.55
$ws.Range("N107").Value = -7065.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 169.8
$ws.Range("I6").Value = 86.75
$ws.Range("J6").Value = 502
$ws.Range("K6").Value = 260.25
$ws.Range("L6").Value = 1506
$ws.Range("M6").Value = -147.25
$ws.Range("N6").Value = -1732
$ws.Range("H29").Value = 143.44444
$ws.Range("I29").Value = 151.71428
$ws.Range("J29").Value = 114.5
$ws.Range("K29").Value = 455.14284
$ws.Range("L29").Value = 343.5
$ws.Range("M29").Value = -178.14284
$ws.Range("N29").Value = -897.5
$ws.Range("H32").Value = 40029.305
$ws.Range("I32").Value = 29999
$ws.Range("J32").Value = 66776.78
$ws.Range("K32").Value = 89997
$ws.Range("L32").Value = 200330.34
$ws.Range("M32").Value = -89714
$ws.Range("N32").Value = -200896.34
$ws.Range("H46").Value = 162165460
$ws.Range("I46").Value = 22428.111
$ws.Range("J46").Value = 370635070
$ws.Range("K46").Value = 67284.333
$ws.Range("L46").Value = 1111905210
$ws.Range("M46").Value = -67193.333
$ws.Range("N46").Value = -1111905392
$ws.Range("H68").Value = 1247.125
$ws.Range("I68").Value = 1057.625
$ws.Range("J68").Value = 1436.625
$ws.Range("K68").Value = 3172.875
$ws.Range("L68").Value = 4309.875
$ws.Range("M68").Value = -2361.875
$ws.Range("N68").Value = -5931.875
$ws.Range("H71").Value = 1247.125
$ws.Range("I71").Value = 1057.625
$ws.Range("J71").Value = 1436.625
$ws.Range("K71").Value = 9518.625
$ws.Range("L71").Value = 12929.625
$ws.Range("M71").Value = -5462.625
$ws.Range("N71").Value = -21041.625
$ws.Range("H86").Value = 223
$ws.Range("I86").Value = 223
$ws.Range("K86").Value = 669
$ws.Range("M86").Value = 517
$ws.Range("H89").Value = 223
$ws.Range("I89").Value = 223
$ws.Range("K89").Value = 2007
$ws.Range("M89").Value = 3921
$ws.Range("H104").Value = 15410
$ws.Range("I104").Value = 8262.6
$ws.Range("K104").Value = 24787.8
$ws.Range("M104").Value = -22166.8
$ws.Range("H113").Value = 2340.889
$ws.Range("I113").Value = 2465.3333
$ws.Range("K113").Value = 7395.999899999999
$ws.Range("M113").Value = -5225.999899999999
$ws.Range("H137").Value = 4549.4443
$ws.Range("I137").Value = 3707.8572
$ws.Range("K137").Value = 11123.5716
$ws.Range("M137").Value = -6023.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3609.1538
$ws.Range("I132").Value = 2901.7878
$ws.Range("J132").Value = 7499.6665
$ws.Range("K132").Value = 8705.3634
$ws.Range("L132").Value = 22498.9995
$ws.Range("M132").Value = -6175.3634
$ws.Range("N132").Value = -27558.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3236.0952
$ws.Range("I7").Value = 2697.875
$ws.Range("K7").Value = 2697.875
$ws.Range("M7").Value = -2585.875
$ws.Range("H122").Value = 4599.8335
$ws.Range("I122").Value = 3519.8
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 10559.4
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -8109.400000000001
$ws.Range("N122").Value = -34900
$ws.Range("H126").Value = 3236.0952
$ws.Range("I126").Value = 2697.875
$ws.Range("K126").Value = 8093.625
$ws.Range("M126").Value = -5623.625
$ws.Range("H136").Value = 10597444
$ws.Range("I136").Value = 12003389
$ws.Range("J136").Value = 52852
$ws.Range("K136").Value = 36010167
$ws.Range("L136").Value = 158556
$ws.Range("M136").Value = -36007617
$ws.Range("N136").Value = -163656

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1778.3846
$ws.Range("I107").Value = 892.8333
$ws.Range("J107").Value = 2537.4285
$ws.Range("K107").Value = 2678.4999
$ws.Range("L107").Value = 7612.2855
$ws.Range("M107").Value = -758.4998999999998
$ws.Range("N107").Value = -11452.2855
$ws.Range("H132").Value = 4665.3696
$ws.Range("I132").Value = 4758.619
$ws.Range("J132").Value = 3686.25
$ws.Range("K132").Value = 14275.857
$ws.Range("L132").Value = 11058.75
$ws.Range("M132").Value = -11745.857
$ws.Range("N132").Value = -16118.75
